$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       GradientBoostingRegressor())]),`n                                            param_grid={'model__max_depth': [3,`n                                                                             5,`n                                                                             7],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# New header cell F1 "Modelo", matching style of other header cells (A1:E1)
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill F2:F5 with the model description text
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText
$ws.Range("F5").Value = $modelText

# Minor floating point precision corrections for B/D columns
$ws.Range("B2").Value = 0.4108830721519875
$ws.Range("D2").Value = 0.525669270734418
$ws.Range("B3").Value = 0.1647362719327807
$ws.Range("D3").Value = 0.321570631223101
$ws.Range("B4").Value = 0.2192987281846222
$ws.Range("D4").Value = 0.3777406617731504
$ws.Range("B5").Value = 0.3638702225807678
$ws.Range("D5").Value = 0.4682092966157618

# Keep rows at their default height - undo any auto row-height expansion
# triggered by entering the long multi-line text in column F.
$ws.Rows.Item(1).EntireRow.AutoFit()
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()
